$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(6,13,14,16,19,25,26,27,32,37,39,43,46,48,50,51,54,55,56,61,62,63,67,69,73,75,76,77,84,86,92,93,95,96,101,105,106,110,116,121,124,125,128,129,132,136,139,143,144,145,149,156,160,161,165,170,171,173,175,177,184,188,191,192)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
